# Scheduled market-data refresh: update cached price figures across the
# Leve profit sheets, and drop the (now stale) priced columns for a block
# of ALC leves that no longer have current pricing data.

$wb = $excel.ActiveWorkbook

# --- ALC: clear currentAveragePrice..LeveProfitHQ (H:N) for rows 125-141 ---
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H125:N141").ClearContents()

# --- ARM: refreshed pricing for rows 3, 108, 110 ---
$wsARM = $wb.Worksheets.Item("ARM")

$wsARM.Range("H3").Value = 3800.9092
$wsARM.Range("I3").Value = 2707.5
$wsARM.Range("J3").Value = 4425.7144
$wsARM.Range("K3").Value = 2707.5
$wsARM.Range("L3").Value = 4425.7144
$wsARM.Range("M3").Value = -2592.5
$wsARM.Range("N3").Value = -4655.7144

$wsARM.Range("H108").Value = 29828
$wsARM.Range("J108").Value = 29828
$wsARM.Range("L108").Value = 29828
$wsARM.Range("N108").Value = -37508

$wsARM.Range("H110").Value = 4426.2354
$wsARM.Range("I110").Value = 4077.875
$wsARM.Range("J110").Value = 10000
$wsARM.Range("K110").Value = 4077.875
$wsARM.Range("L110").Value = 10000
$wsARM.Range("M110").Value = -2032.875
$wsARM.Range("N110").Value = -14090

# --- CRP: refreshed pricing for row 123 ---
$wsCRP = $wb.Worksheets.Item("CRP")

$wsCRP.Range("H123").Value = 50592
$wsCRP.Range("J123").Value = 50592
$wsCRP.Range("L123").Value = 50592
$wsCRP.Range("N123").Value = -60392

# --- CUL: refreshed pricing for rows 131, 138 ---
$wsCUL = $wb.Worksheets.Item("CUL")

$wsCUL.Range("H131").Value = 838.25
$wsCUL.Range("I131").Value = 330.83334
$wsCUL.Range("J131").Value = 927.7941
$wsCUL.Range("K131").Value = 992.5000200000001
$wsCUL.Range("L131").Value = 2783.3823
$wsCUL.Range("M131").Value = 4047.49998
$wsCUL.Range("N131").Value = -12863.3823

$wsCUL.Range("H138").Value = 45456280
$wsCUL.Range("I138").Value = 76924200
$wsCUL.Range("J138").Value = 2611.111
$wsCUL.Range("K138").Value = 230772600
$wsCUL.Range("L138").Value = 7833.333
$wsCUL.Range("M138").Value = -230767460
$wsCUL.Range("N138").Value = -18113.333

# --- GSM: refreshed pricing for rows 70, 73, 108 ---
$wsGSM = $wb.Worksheets.Item("GSM")

$wsGSM.Range("H70").Value = 20377.88
$wsGSM.Range("I70").Value = 43100
$wsGSM.Range("J70").Value = 5040.45
$wsGSM.Range("K70").Value = 43100
$wsGSM.Range("L70").Value = 5040.45
$wsGSM.Range("M70").Value = -42830
$wsGSM.Range("N70").Value = -5580.45

$wsGSM.Range("H73").Value = 20377.88
$wsGSM.Range("I73").Value = 43100
$wsGSM.Range("J73").Value = 5040.45
$wsGSM.Range("K73").Value = 43100
$wsGSM.Range("L73").Value = 5040.45
$wsGSM.Range("M73").Value = -42164
$wsGSM.Range("N73").Value = -6912.45

$wsGSM.Range("H108").Value = 39888
$wsGSM.Range("J108").Value = 39888
$wsGSM.Range("L108").Value = 39888
$wsGSM.Range("N108").Value = -47568
